# Update the active selection on Sheet1 (Decision Table) from A20 to A15,
# and turn on the "Print Gridlines" option for that sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Decision Table")
$ws.Activate()
$ws.Range("A15").Select() | Out-Null
$ws.PageSetup.PrintGridlines = $true
